$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.947.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.367.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.30"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.365.71"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.337"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.13"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.788.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.799.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.338.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.77"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.13%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +8.84%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "316.55"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.59"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.73%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.35%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.27%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +15.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.385"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.11"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "318.87"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.30"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.40"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0214"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.97%  "
